$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 197
$ws.Range("E2").Value = "'123,852,103.00"
$ws.Range("G2").Value = "'22.94"
$ws.Range("H2").Value = 33.354999999999997
$ws.Range("J2").Value = 117
$ws.Range("K2").Value = "'15.15"
$ws.Range("L2").Value = "'15.29"

# Row 3
$ws.Range("D3").Value = 187
$ws.Range("E3").Value = "'115,653,756.00"
$ws.Range("G3").Value = "'19.97"
$ws.Range("H3").Value = 30.991
$ws.Range("J3").Value = 111
$ws.Range("K3").Value = "'13.43"
$ws.Range("L3").Value = "'13.94"

# Row 4
$ws.Range("D4").Value = 224
$ws.Range("E4").Value = "'144,865,306.00"
$ws.Range("G4").Value = "'25.50"
$ws.Range("H4").Value = 30.26
$ws.Range("J4").Value = 123
$ws.Range("K4").Value = "'15.95"
$ws.Range("L4").Value = "'15.49"

# Row 5
$ws.Range("D5").Value = 164
$ws.Range("E5").Value = "'102,001,161.00"
$ws.Range("G5").Value = "'19.00"
$ws.Range("H5").Value = 18.739999999999998
$ws.Range("J5").Value = 98
$ws.Range("K5").Value = "'13.39"
$ws.Range("L5").Value = "'12.76"

# Row 6
$ws.Range("D6").Value = 182
$ws.Range("E6").Value = "'118,228,866.00"
$ws.Range("G6").Value = "'19.94"
$ws.Range("H6").Value = 34.418999999999997
$ws.Range("J6").Value = 96
$ws.Range("K6").Value = "'12.45"
$ws.Range("L6").Value = "'12.09"

# Row 7
$ws.Range("D7").Value = 196
$ws.Range("E7").Value = "'151,065,305.00"
$ws.Range("G7").Value = "'26.08"
$ws.Range("H7").Value = 32.862000000000002
$ws.Range("J7").Value = 99
$ws.Range("K7").Value = "'13.52"
$ws.Range("L7").Value = "'12.41"

# Row 8
$ws.Range("D8").Value = 201
$ws.Range("E8").Value = "'134,635,427.00"
$ws.Range("G8").Value = "'21.87"
$ws.Range("H8").Value = 57.741
$ws.Range("J8").Value = 119
$ws.Range("K8").Value = "'13.84"
$ws.Range("L8").Value = "'14.89"

# Row 9
$ws.Range("D9").Value = 169
$ws.Range("E9").Value = "'110,136,941.00"
$ws.Range("G9").Value = "'18.69"
$ws.Range("H9").Value = 34.546999999999997
$ws.Range("J9").Value = 86
$ws.Range("K9").Value = "'10.78"
$ws.Range("L9").Value = "'10.75"

# Row 10
$ws.Range("D10").Value = 243
$ws.Range("E10").Value = "'146,746,736.00"
$ws.Range("G10").Value = "'25.31"
$ws.Range("H10").Value = 43.19
$ws.Range("J10").Value = 131
$ws.Range("K10").Value = "'14.35"
$ws.Range("L10").Value = "'16.46"

# Row 11
$ws.Range("D11").Value = 207
$ws.Range("E11").Value = "'146,067,991.00"
$ws.Range("G11").Value = "'24.55"
$ws.Range("H11").Value = 40.686999999999998
$ws.Range("J11").Value = 126
$ws.Range("K11").Value = "'14.29"
$ws.Range("L11").Value = "'15.85"

# Row 12
$ws.Range("D12").Value = 182
$ws.Range("E12").Value = "'121,060,527.00"
$ws.Range("G12").Value = "'20.27"
$ws.Range("H12").Value = 37.368000000000002
$ws.Range("J12").Value = 110
$ws.Range("K12").Value = "'13.20"
$ws.Range("L12").Value = "'13.80"

# Row 13
$ws.Range("D13").Value = 188
$ws.Range("E13").Value = "'125,642,039.00"
$ws.Range("G13").Value = "'21.56"
$ws.Range("H13").Value = 37.497
$ws.Range("J13").Value = 109
$ws.Range("K13").Value = "'14.38"
$ws.Range("L13").Value = "'13.57"

# Row 14
$ws.Range("D14").Value = 226
$ws.Range("E14").Value = "'129,742,447.00"
$ws.Range("G14").Value = "'22.74"
$ws.Range("H14").Value = 44.610999999999997
$ws.Range("J14").Value = 107
$ws.Range("K14").Value = "'11.52"
$ws.Range("L14").Value = "'13.33"

# Row 15
$ws.Range("D15").Value = 177
$ws.Range("E15").Value = "'148,003,578.00"
$ws.Range("G15").Value = "'25.68"
$ws.Range("H15").Value = 61.04
$ws.Range("J15").Value = 88
$ws.Range("K15").Value = "'13.78"
$ws.Range("L15").Value = "'11.06"
